$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: add the new trailing header columns (E1:M1), reusing the same
# header style (bold, centered, bordered) already applied to A1:D1. A plain
# Value assignment won't carry the style, so copy A1's formatting and paste
# it (formats only) onto each new header cell after setting its text.
$headers = @(
    @("E1", "dose_amt"),
    @("F1", "nda_num"),
    @("G1", "route"),
    @("H1", "dose_unit"),
    @("I1", "dose_form"),
    @("J1", "dose_freq"),
    @("K1", "dechal"),
    @("L1", "rechal"),
    @("M1", "role_cod")
)

$ws.Range("A1").Copy()
foreach ($pair in $headers) {
    $addr = $pair[0]
    $text = $pair[1]
    $ws.Range($addr).Value = $text
    $ws.Range($addr).PasteSpecial(-4122)  # xlPasteFormats
}

# --- Row 2: fill in default/sample values.
# A2 and B2 remain blank (unchanged). C2/D2 get text defaults; E2/F2 are
# numeric defaults; G2:L2 get text defaults; M2 gets its own text value.
$ws.Range("C2").Value = "Unknown"
$ws.Range("D2").Value = "Unknown"
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = "Unknown"
$ws.Range("H2").Value = "Unknown"
$ws.Range("I2").Value = "Unknown"
$ws.Range("J2").Value = "Unknown"
$ws.Range("K2").Value = "Unknown"
$ws.Range("L2").Value = "Unknown"
$ws.Range("M2").Value = "PS"
